$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 40000200
$ws.Range("B6").Value = "AreaAttack/Prefabs/AreaAttack_Ice"
$ws.Range("C6").Value = "Magic"
$ws.Range("D6").Value = 0

$ws.Range("G11").Select() | Out-Null
